$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 20) to the finStatID/guestListID/talaID register.
$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 1).Style = $ws.Cells.Item(19, 1).Style

$ws.Cells.Item(20, 2).Value = "G2021092002"
$ws.Cells.Item(20, 3).Value = "F2021092001"
$ws.Cells.Item(20, 4).Value = "1"
